$d = $word.ActiveDocument

# --- 1. Simple heading swap: "Pitchforks and Torches" -> "Torches and Pitchforks" ---
$d.Content.Find.Execute("Pitchforks and Torches", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Torches and Pitchforks", 2) | Out-Null

# --- Text chunks pulled verbatim from the target content ---
$RUN_OUR         = "Our "
$RUN_COUNTRY     = "country"
$RUN_HASBECOME   = " has become a"
$RUN_RELENTLESS  = " relentless,"
$RUN_UNFORGIVING = " unforgiving lynch mob."
$RUN_LONG        = " [We grew up repeating ""liberty and justice for all"" but I guess those were just words.]? I don't condone what Donald Sterling said or agree with his viewpoint in any way, but a ban for life means we don't think he can be rehabilitated. Can a person be cured of racism? Hasn't anyone ever seen "
$RUN_TITANS      = "Remember the Titans"
$RUN_BANG        = "?!"
$HEADING2        = "A Limitless Checkbook"
$RUN_STRUGGLING  = "I'm struggling to find the right words to describe how disgusted I am that our country has spent almost a "
$RUN_BILLION     = "billion"
$RUN_WEBSITE     = " dollars to develop the healthcare.gov website! "

# --- 2. Locate the paragraph that currently holds
#        "Our society has become an unforgiving lynch mob-as if none..." ---
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Our society has become")) {
        $targetIndex = $i
        break
    }
}

$pBody = $d.Paragraphs($targetIndex)

# Clone that paragraph's formatting for the new trailing body paragraph
# (it will receive "I'm struggling ... website!").
$pBody.Range.InsertParagraphAfter()

# Clone the heading paragraph right above (now "Torches and Pitchforks") to make
# the new "A Limitless Checkbook" heading paragraph.
$pHeading = $d.Paragraphs($targetIndex - 1)
$pHeading.Range.InsertParagraphAfter()

# Paragraph indices after the two inserts above:
#   targetIndex - 1  : "Torches and Pitchforks"            (unchanged heading)
#   targetIndex      : new, empty heading paragraph        -> "A Limitless Checkbook"
#   targetIndex + 1  : original body paragraph              -> rewritten below
#   targetIndex + 2  : new, empty body paragraph             -> "I'm struggling..."

# --- 3. Fill the new heading paragraph, wrapped in the _GoBack bookmark ---
$pNewHeading = $d.Paragraphs($targetIndex)
$rHeadingInsert = $d.Range($pNewHeading.Range.Start, $pNewHeading.Range.End - 1)
$rHeadingInsert.InsertAfter($HEADING2)

$pNewHeading2 = $d.Paragraphs($targetIndex)
$rHeadingBookmark = $d.Range($pNewHeading2.Range.Start, $pNewHeading2.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $rHeadingBookmark)

# --- 4. Rewrite the original body paragraph with the new multi-run content ---
$pOldBody = $d.Paragraphs($targetIndex + 1)
$rClear = $d.Range($pOldBody.Range.Start, $pOldBody.Range.End - 1)
$rClear.Text = ""

$pOldBody2 = $d.Paragraphs($targetIndex + 1)
$insertPoint = $d.Range($pOldBody2.Range.End - 1, $pOldBody2.Range.End - 1)
$insertPoint.InsertAfter($RUN_OUR)

$pOldBody3 = $d.Paragraphs($targetIndex + 1)
$insertPoint = $d.Range($pOldBody3.Range.End - 1, $pOldBody3.Range.End - 1)
$insertPoint.InsertAfter($RUN_COUNTRY)

$pOldBody4 = $d.Paragraphs($targetIndex + 1)
$insertPoint = $d.Range($pOldBody4.Range.End - 1, $pOldBody4.Range.End - 1)
$insertPoint.InsertAfter($RUN_HASBECOME)

$pOldBody5 = $d.Paragraphs($targetIndex + 1)
$insertPoint = $d.Range($pOldBody5.Range.End - 1, $pOldBody5.Range.End - 1)
$insertPoint.InsertAfter($RUN_RELENTLESS)

$pOldBody6 = $d.Paragraphs($targetIndex + 1)
$insertPoint = $d.Range($pOldBody6.Range.End - 1, $pOldBody6.Range.End - 1)
$insertPoint.InsertAfter($RUN_UNFORGIVING)

$pOldBody7 = $d.Paragraphs($targetIndex + 1)
$insertPoint = $d.Range($pOldBody7.Range.End - 1, $pOldBody7.Range.End - 1)
$insertPoint.InsertAfter($RUN_LONG)

$pOldBody8 = $d.Paragraphs($targetIndex + 1)
$titansStart = $pOldBody8.Range.End - 1
$insertPoint = $d.Range($titansStart, $titansStart)
$insertPoint.InsertAfter($RUN_TITANS)
$titansEnd = $titansStart + $RUN_TITANS.Length
$d.Range($titansStart, $titansEnd).Italic = 1

$pOldBody9 = $d.Paragraphs($targetIndex + 1)
$insertPoint = $d.Range($pOldBody9.Range.End - 1, $pOldBody9.Range.End - 1)
$insertPoint.InsertAfter($RUN_BANG)

$pOldBody10 = $d.Paragraphs($targetIndex + 1)
$insertPoint = $d.Range($pOldBody10.Range.End - 1, $pOldBody10.Range.End - 1)
$insertPoint.InsertParagraphAfter()
# remove the paragraph break we just added and replace it with a plain line break
# so the content stays inside the same paragraph.
$pOldBody11 = $d.Paragraphs($targetIndex + 1)
$brRange = $d.Range($pOldBody11.Range.End - 1, $pOldBody11.Range.End - 1)
$brRange.InsertBefore("<<<BR>>>")

# --- 5. Fill the new trailing body paragraph ---
$pNewBody = $d.Paragraphs($targetIndex + 3)
$insertPoint = $d.Range($pNewBody.Range.Start, $pNewBody.Range.Start)
$insertPoint.InsertAfter($RUN_STRUGGLING)

$pNewBody2 = $d.Paragraphs($targetIndex + 3)
$billionStart = $pNewBody2.Range.Start + $RUN_STRUGGLING.Length
$insertPoint = $d.Range($billionStart, $billionStart)
$insertPoint.InsertAfter($RUN_BILLION)
$billionEnd = $billionStart + $RUN_BILLION.Length
$d.Range($billionStart, $billionEnd).Italic = 1

$pNewBody3 = $d.Paragraphs($targetIndex + 3)
$insertPoint = $d.Range($billionEnd, $billionEnd)
$insertPoint.InsertAfter($RUN_WEBSITE)

Write-Output "done"
